$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AY1").Value = 0.98100799277227846
$ws.Range("AZ1").Value = 0.60935254108427206
$ws.Range("X2").Value = 0.99889491068425218
$ws.Range("I3").Value = 0.87952289069647205
$ws.Range("U3").Value = 0.78995708406664744
$ws.Range("V3").Value = 0.95695139180798172
$ws.Range("BG3").Value = 0.9223444695591605
$ws.Range("AA4").Value = 0.72271691930243698
$ws.Range("AH4").Value = 0.92609362951465102
$ws.Range("AF5").Value = 0.99599330126479091
$ws.Range("AD6").Value = 0.89809512901822797
$ws.Range("BA6").Value = 0.94894427696226713
$ws.Range("AQ7").Value = 0.93741573282550417
$ws.Range("BI7").Value = 0.79302046123372039
$ws.Range("B8").Value = 0.98866742043488387
$ws.Range("F10").Value = 0.89218918946013315
$ws.Range("AB10").Value = 0.81491884148150961
$ws.Range("J11").Value = 0.8533170113966343
$ws.Range("BB11").Value = 0.96794378491167521
$ws.Range("BM11").Value = 0.80859175277462136
$ws.Range("I12").Value = 0.5890466313842182
$ws.Range("J12").Value = 0.95180313349295487
$ws.Range("AI12").Value = 0.7079044096489886
$ws.Range("BC12").Value = 0.94893471914689553
$ws.Range("B13").Value = 0.75790312414866867
$ws.Range("O13").Value = 0.81696582929235317
$ws.Range("U13").Value = 0.52907748120223286
$ws.Range("Z13").Value = 0.8517583890775775
$ws.Range("BB13").Value = 0.64758492193120376
$ws.Range("AV14").Value = 0.72629590881244943
$ws.Range("R15").Value = 0.94020839920002841
$ws.Range("E16").Value = 0.87370362260992951
$ws.Range("AP16").Value = 0.96641854285614259
$ws.Range("BC17").Value = 0.98989785735960889
$ws.Range("T18").Value = 0.71092829120980283
$ws.Range("BN18").Value = 0.67370499871735334
$ws.Range("R19").Value = 0.8566475988909692
$ws.Range("AY20").Value = 0.78787244405374657
$ws.Range("N21").Value = 0.850780087249428
$ws.Range("T21").Value = 0.95980434877899001
$ws.Range("AB21").Value = 0.95273110332827771
$ws.Range("AJ21").Value = 0.8926583479135014
$ws.Range("AU21").Value = 0.85951016791266754
$ws.Range("X23").Value = 0.61069274801928164
$ws.Range("AU23").Value = 0.89476181646595743
$ws.Range("AA24").Value = 0.89949725510222434
$ws.Range("K25").Value = 0.85615134769784951
$ws.Range("AL25").Value = 0.57534797062848986
$ws.Range("T26").Value = 0.95866711372495905
$ws.Range("V26").Value = 0.84035312651736704
$ws.Range("AE26").Value = 0.69790161714338539
$ws.Range("AG26").Value = 0.94690915404317733
$ws.Range("K27").Value = 0.96988331952710338
$ws.Range("BM27").Value = 0.72336209224883341
$ws.Range("E28").Value = 0.74423908131630756
$ws.Range("C29").Value = 0.58571292788945506
$ws.Range("G29").Value = 0.63416178958682323
$ws.Range("Q29").Value = 0.65675635802740162
$ws.Range("AA29").Value = 0.94586637267422513
$ws.Range("BE29").Value = 0.71416469932074522
$ws.Range("V30").Value = 0.99527958205171063
$ws.Range("AX31").Value = 0.85619141031553658
$ws.Range("H32").Value = 0.87487384376903199
$ws.Range("S32").Value = 0.97560575494342705
$ws.Range("AJ32").Value = 0.98779170890950807
$ws.Range("AN32").Value = 0.85446407120367129
$ws.Range("Q33").Value = 0.95029014668267497
$ws.Range("BD34").Value = 0.6598351122856585
$ws.Range("AA35").Value = 0.97146059160435971
$ws.Range("AW35").Value = 0.93438941217991478
$ws.Range("K36").Value = 0.93569759700572819
$ws.Range("AX36").Value = 0.93462319682772044
$ws.Range("AL37").Value = 0.95276442407505357
$ws.Range("AM37").Value = 0.94684249892323802
$ws.Range("H39").Value = 0.75430156206977128
$ws.Range("Q39").Value = 0.77042656775423524
$ws.Range("AA39").Value = 0.99394970621820122
$ws.Range("AW40").Value = 0.62275874863835567
$ws.Range("AD41").Value = 0.88606820189602264
$ws.Range("BC41").Value = 0.91251874162281488
$ws.Range("AL42").Value = 0.9497874598213667
$ws.Range("AR43").Value = 0.79358415958761186
$ws.Range("BH43").Value = 0.95271434998001081
$ws.Range("J44").Value = 0.92255691340705981
$ws.Range("AG45").Value = 0.51579298917022398
$ws.Range("H46").Value = 0.83243935132915969
$ws.Range("AZ46").Value = 0.98482829200200261
$ws.Range("BM46").Value = 0.70081569195637994
$ws.Range("S47").Value = 0.92267873045500193
$ws.Range("X47").Value = 0.56822626670333598
$ws.Range("O49").Value = 0.89342998409985663
$ws.Range("Y49").Value = 0.96240729256836643
$ws.Range("AL49").Value = 0.73155056848378541
$ws.Range("Z50").Value = 0.83717592336834712
$ws.Range("AU50").Value = 0.71403649215722687
$ws.Range("BP51").Value = 0.78992641880897896
$ws.Range("H52").Value = 0.79936663304003064
$ws.Range("S52").Value = 0.7811934422113922
$ws.Range("AC52").Value = 0.90126574553090744
$ws.Range("AS52").Value = 0.99589378122277117
$ws.Range("BA54").Value = 0.79908329906492903
$ws.Range("BA55").Value = 0.81203461937473898
$ws.Range("BK55").Value = 0.94208164504232061
$ws.Range("BL55").Value = 0.64478970233008814
$ws.Range("H56").Value = 0.74073456049553021
$ws.Range("C57").Value = 0.8802620143771489
$ws.Range("AR57").Value = 0.81177606376313816
$ws.Range("BG57").Value = 0.71215891083941596
$ws.Range("AT58").Value = 0.96097852082959834
$ws.Range("BD58").Value = 0.81655355333348512
$ws.Range("BH58").Value = 0.87053698111233302
$ws.Range("P59").Value = 0.75885549024834387
$ws.Range("AV59").Value = 0.91477585284354168
$ws.Range("BF59").Value = 0.72317762756089654
$ws.Range("H60").Value = 0.95596229295260127
$ws.Range("O60").Value = 0.61346668644432367
$ws.Range("AX60").Value = 0.63292237656309713
$ws.Range("BP60").Value = 0.96197147284242124
$ws.Range("AK61").Value = 0.90990936217959906
$ws.Range("AN61").Value = 0.81980113646151098
$ws.Range("BG61").Value = 0.98351135613209806
$ws.Range("BK62").Value = 0.73389190644148683
$ws.Range("AX63").Value = 0.59606120567357479
$ws.Range("BM63").Value = 0.93367298497943008
$ws.Range("AJ64").Value = 0.82184063916314121
$ws.Range("BJ64").Value = 0.80881060938940919
$ws.Range("BN64").Value = 0.89714484046777021
$ws.Range("M66").Value = 0.99136761734219359
$ws.Range("BK66").Value = 0.99726060124088001
$ws.Range("A67").Value = 0.97308299850772439
$ws.Range("AS67").Value = 0.80840348604343326
$ws.Range("AW67").Value = 0.90472496254484991
$ws.Range("I68").Value = 0.76845952830329023
$ws.Range("AH68").Value = 0.88666331915756835
$ws.Range("AP68").Value = 0.97406482570285891
$ws.Range("AV68").Value = 0.99369965776871561
